$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 914
$ws.Range("I38").Value = 233.16667
$ws.Range("K38").Value = 699.50001
$ws.Range("M38").Value = -327.50001
$ws.Range("H41").Value = 1022.6818
$ws.Range("J41").Value = 916.8333
$ws.Range("L41").Value = 916.8333
$ws.Range("N41").Value = -1796.8333
$ws.Range("H53").Value = 73.833336
$ws.Range("J53").Value = 75.59999999999999
$ws.Range("L53").Value = 75.59999999999999
$ws.Range("N53").Value = -1349.6
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("H87").Value = 249997.5
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 249997.5
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H94").Value = 653.1818
$ws.Range("I94").Value = 576.3333
$ws.Range("K94").Value = 576.3333
$ws.Range("M94").Value = -125.3333
$ws.Range("H98").Value = 2203.2144
$ws.Range("I98").Value = 570.6667
$ws.Range("J98").Value = 11998.5
$ws.Range("K98").Value = 570.6667
$ws.Range("L98").Value = 11998.5
$ws.Range("M98").Value = 927.3333
$ws.Range("N98").Value = -14994.5
$ws.Range("H100").Value = 580
$ws.Range("I100").Value = 580
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 580
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -39
$ws.Range("H106").Value = 5886.263
$ws.Range("I106").Value = 5886.263
$ws.Range("K106").Value = 5886.263
$ws.Range("M106").Value = -5255.263
$ws.Range("H112").Value = 2838.8076
$ws.Range("J112").Value = 2929.75
$ws.Range("L112").Value = 8789.25
$ws.Range("N112").Value = -11005.25
$ws.Range("H113").Value = 2113.25
$ws.Range("I113").Value = 2317.6667
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2317.6667
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 936.3332999999998
$ws.Range("N113").Value = -8008
$ws.Range("H122").Value = 2203.2144
$ws.Range("I122").Value = 570.6667
$ws.Range("J122").Value = 11998.5
$ws.Range("K122").Value = 1712.0001
$ws.Range("L122").Value = 35995.5
$ws.Range("M122").Value = 737.9999
$ws.Range("N122").Value = -40895.5
$ws.Range("H125").Value = 4845.3335
$ws.Range("I125").Value = 5374.9
$ws.Range("J125").Value = 2197.5
$ws.Range("K125").Value = 48374.1
$ws.Range("L125").Value = 19777.5
$ws.Range("M125").Value = -45914.1
$ws.Range("N125").Value = -24697.5
$ws.Range("H132").Value = 4106.5
$ws.Range("I132").Value = 3711.0688
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 11133.2064
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -8603.206399999999
$ws.Range("N132").Value = -24260
$ws.Range("H135").Value = 4816.3887
$ws.Range("I135").Value = 4779.6
$ws.Range("J135").Value = 4862.375
$ws.Range("K135").Value = 43016.4
$ws.Range("L135").Value = 43761.375
$ws.Range("M135").Value = -40481.4
$ws.Range("N135").Value = -48831.375
$ws.Range("H137").Value = 1025.0541
$ws.Range("I137").Value = 997.71875
$ws.Range("K137").Value = 2993.15625
$ws.Range("M137").Value = -443.15625
$ws.Range("H138").Value = 3621.1965
$ws.Range("I138").Value = 1383.6
$ws.Range("J138").Value = 4439.829
$ws.Range("K138").Value = 4150.799999999999
$ws.Range("L138").Value = 13319.487
$ws.Range("M138").Value = 989.2000000000007
$ws.Range("N138").Value = -23599.487
$ws.Range("H141").Value = 2013.7142
$ws.Range("I141").Value = 2013.7142
$ws.Range("K141").Value = 6041.142599999999
$ws.Range("M141").Value = -861.1425999999992
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8578.200000000001
$ws.Range("I2").Value = 8578.200000000001
$ws.Range("K2").Value = 8578.200000000001
$ws.Range("M2").Value = -8465.200000000001
$ws.Range("H4").Value = 349.7647
$ws.Range("I4").Value = 282.375
$ws.Range("J4").Value = 409.66666
$ws.Range("K4").Value = 282.375
$ws.Range("L4").Value = 409.66666
$ws.Range("M4").Value = -166.375
$ws.Range("N4").Value = -641.66666
$ws.Range("H32").Value = 7966.5454
$ws.Range("I32").Value = 2395.92
$ws.Range("K32").Value = 2395.92
$ws.Range("M32").Value = -2108.92
$ws.Range("H61").Value = 4049.1316
$ws.Range("I61").Value = 2726.1785
$ws.Range("J61").Value = 7753.4
$ws.Range("K61").Value = 2726.1785
$ws.Range("L61").Value = 7753.4
$ws.Range("M61").Value = -2514.1785
$ws.Range("N61").Value = -8177.4
$ws.Range("H76").Value = 43532.375
$ws.Range("J76").Value = 43532.375
$ws.Range("L76").Value = 43532.375
$ws.Range("N76").Value = -44208.375
$ws.Range("H79").Value = 43532.375
$ws.Range("J79").Value = 43532.375
$ws.Range("L79").Value = 43532.375
$ws.Range("N79").Value = -45872.375
$ws.Range("H116").Value = 8578.200000000001
$ws.Range("I116").Value = 8578.200000000001
$ws.Range("K116").Value = 8578.200000000001
$ws.Range("M116").Value = -6284.200000000001
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0
$ws.Range("H122").Value = 4626.087
$ws.Range("I122").Value = 3438.8125
$ws.Range("K122").Value = 10316.4375
$ws.Range("M122").Value = -7866.4375
$ws.Range("H132").Value = 4203.8
$ws.Range("I132").Value = 3975.8333
$ws.Range("K132").Value = 11927.4999
$ws.Range("M132").Value = -9397.499899999999
$ws.Range("H134").Value = 89664
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 89664
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").Value = 89664
$ws.Range("N134").Value = -99804
$ws.Range("H136").Value = 4049.1316
$ws.Range("I136").Value = 2726.1785
$ws.Range("J136").Value = 7753.4
$ws.Range("K136").Value = 8178.5355
$ws.Range("L136").Value = 23260.2
$ws.Range("M136").Value = -5628.5355
$ws.Range("N136").Value = -28360.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8578.200000000001
$ws.Range("I3").Value = 8578.200000000001
$ws.Range("K3").Value = 8578.200000000001
$ws.Range("M3").Value = -8464.200000000001
$ws.Range("H76").Value = 28337.059
$ws.Range("J76").Value = 28545.688
$ws.Range("L76").Value = 28545.688
$ws.Range("N76").Value = -29175.688
$ws.Range("H79").Value = 28337.059
$ws.Range("J79").Value = 28545.688
$ws.Range("L79").Value = 28545.688
$ws.Range("N79").Value = -30729.688
$ws.Range("H87").Value = 109994.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 109994.5
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("M87").Value = 109994.5
$ws.Range("N87").Value = -112490.5
$ws.Range("H90").Value = 109994.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 109994.5
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("M90").Value = 329983.5
$ws.Range("N90").Value = -342463.5
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 929.2857
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 1000
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 4158.615
$ws.Range("I134").Value = 4076.5454
$ws.Range("K134").Value = 12229.6362
$ws.Range("M134").Value = -9694.636200000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2335
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = 3500
$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = 135
$ws.Range("N11").Value = -3780
$ws.Range("H22").Value = 1788.1666
$ws.Range("J22").Value = 3500
$ws.Range("L22").Value = 3500
$ws.Range("N22").Value = -4200
$ws.Range("H58").Value = 2575.9333
$ws.Range("I58").Value = 2923.5
$ws.Range("J58").Value = 2178.7144
$ws.Range("K58").Value = 2923.5
$ws.Range("L58").Value = 2178.7144
$ws.Range("M58").Value = -2720.5
$ws.Range("N58").Value = -2584.7144
$ws.Range("H93").Value = 21000
$ws.Range("I93").Value = 21000
$ws.Range("K93").Value = 21000
$ws.Range("M93").Value = -19128
$ws.Range("H122").Value = 90698.266
$ws.Range("I122").Value = 160788.95
$ws.Range("J122").Value = 1916.7333
$ws.Range("K122").Value = 482366.85
$ws.Range("L122").Value = 5750.199900000001
$ws.Range("M122").Value = -479916.85
$ws.Range("N122").Value = -10650.1999
$ws.Range("H132").Value = 1547.1818
$ws.Range("I132").Value = 1457.579
$ws.Range("K132").Value = 4372.737
$ws.Range("M132").Value = -1842.737
$ws.Range("H134").Value = 1387.5555
$ws.Range("I134").Value = 997.5
$ws.Range("J134").Value = 1699.6
$ws.Range("K134").Value = 2992.5
$ws.Range("L134").Value = 5098.799999999999
$ws.Range("M134").Value = -457.5
$ws.Range("N134").Value = -10168.8
$ws.Range("H136").Value = 2575.9333
$ws.Range("I136").Value = 2923.5
$ws.Range("J136").Value = 2178.7144
$ws.Range("K136").Value = 8770.5
$ws.Range("L136").Value = 6536.1432
$ws.Range("M136").Value = -6220.5
$ws.Range("N136").Value = -11636.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 829.95
$ws.Range("I2").Value = 1082.2142
$ws.Range("J2").Value = 241.33333
$ws.Range("K2").Value = 6493.285199999999
$ws.Range("L2").Value = 1447.99998
$ws.Range("M2").Value = -6380.285199999999
$ws.Range("N2").Value = -1673.99998
$ws.Range("H5").Value = 405.03705
$ws.Range("I5").Value = 405.03705
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1215.11115
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1103.11115
$ws.Range("H39").Value = 6096.154
$ws.Range("I39").Value = 3044.875
$ws.Range("J39").Value = 10978.2
$ws.Range("K39").Value = 9134.625
$ws.Range("L39").Value = 32934.60000000001
$ws.Range("M39").Value = -8840.625
$ws.Range("N39").Value = -33522.60000000001
$ws.Range("H92").Value = 727
$ws.Range("I92").Value = 600
$ws.Range("K92").Value = 1800
$ws.Range("M92").Value = -552
$ws.Range("H135").Value = 405.03705
$ws.Range("I135").Value = 405.03705
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3645.33345
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -1110.33345
$ws.Range("H136").Value = 6072.4443
$ws.Range("I136").Value = 4808.353
$ws.Range("K136").Value = 14425.059
$ws.Range("M136").Value = -9325.059000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5133.375
$ws.Range("I80").Value = 3700.5
$ws.Range("K80").Value = 3700.5
$ws.Range("M80").Value = -2702.5
$ws.Range("H83").Value = 5133.375
$ws.Range("I83").Value = 3700.5
$ws.Range("K83").Value = 18502.5
$ws.Range("M83").Value = -13510.5
$ws.Range("H102").Value = 2010.5333
$ws.Range("I102").Value = 1958.3077
$ws.Range("K102").Value = 1958.3077
$ws.Range("M102").Value = -336.3077000000001
$ws.Range("H122").Value = 5133.364
$ws.Range("I122").Value = 4940.3335
$ws.Range("K122").Value = 14821.0005
$ws.Range("M122").Value = -12371.0005
$ws.Range("H126").Value = 9355.143
$ws.Range("J126").Value = 12499.5
$ws.Range("L126").Value = 37498.5
$ws.Range("N126").Value = -42438.5
$ws.Range("H132").Value = 3939.9473
$ws.Range("I132").Value = 4000.0212
$ws.Range("J132").Value = 3657.6
$ws.Range("K132").Value = 12000.0636
$ws.Range("L132").Value = 10972.8
$ws.Range("M132").Value = -9470.063600000001
$ws.Range("N132").Value = -16032.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3358
$ws.Range("I7").Value = 3271.6316
$ws.Range("K7").Value = 3271.6316
$ws.Range("M7").Value = -3159.6316
$ws.Range("H40").Value = 3089.7778
$ws.Range("I40").Value = 2592.1667
$ws.Range("J40").Value = 4085
$ws.Range("K40").Value = 2592.1667
$ws.Range("L40").Value = 4085
$ws.Range("M40").Value = -2456.1667
$ws.Range("N40").Value = -4357
$ws.Range("H46").Value = 4570
$ws.Range("I46").Value = 2874.75
$ws.Range("J46").Value = 6830.3335
$ws.Range("K46").Value = 2874.75
$ws.Range("L46").Value = 6830.3335
$ws.Range("M46").Value = -2686.75
$ws.Range("N46").Value = -7206.3335
$ws.Range("H100").Value = 2992.25
$ws.Range("I100").Value = 1597.2
$ws.Range("K100").Value = 1597.2
$ws.Range("M100").Value = -1056.2
$ws.Range("H122").Value = 6910.9546
$ws.Range("I122").Value = 7002.4287
$ws.Range("K122").Value = 21007.2861
$ws.Range("M122").Value = -18557.2861
$ws.Range("H126").Value = 3358
$ws.Range("I126").Value = 3271.6316
$ws.Range("K126").Value = 9814.8948
$ws.Range("M126").Value = -7344.8948
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 4099.8335
$ws.Range("I136").Value = 3225.1428
$ws.Range("J136").Value = 7161.25
$ws.Range("K136").Value = 9675.428400000001
$ws.Range("L136").Value = 21483.75
$ws.Range("M136").Value = -7125.428400000001
$ws.Range("N136").Value = -26583.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 200007
$ws.Range("J18").Value = 200007
$ws.Range("L18").Value = 200007
$ws.Range("N18").Value = -200353
$ws.Range("H81").Value = 2642.5386
$ws.Range("I81").Value = 2724.5
$ws.Range("K81").Value = 5449
$ws.Range("M81").Value = -4388
$ws.Range("H84").Value = 2642.5386
$ws.Range("I84").Value = 2724.5
$ws.Range("K84").Value = 27245
$ws.Range("M84").Value = -21941
$ws.Range("H122").Value = 2577.158
$ws.Range("I122").Value = 2919.2
$ws.Range("J122").Value = 2197.111
$ws.Range("K122").Value = 8757.599999999999
$ws.Range("L122").Value = 6591.333
$ws.Range("M122").Value = -6307.599999999999
$ws.Range("N122").Value = -11491.333
$ws.Range("H126").Value = 2082
$ws.Range("I126").Value = 2178.5
$ws.Range("J126").Value = 924
$ws.Range("K126").Value = 6535.5
$ws.Range("L126").Value = 2772
$ws.Range("M126").Value = -4065.5
$ws.Range("N126").Value = -7712
$ws.Range("H132").Value = 8720.777
$ws.Range("I132").Value = 6415.4165
$ws.Range("K132").Value = 19246.2495
$ws.Range("M132").Value = -16716.2495
$ws.Range("H136").Value = 3842.2942
$ws.Range("I136").Value = 3653.0344
$ws.Range("J136").Value = 4940
$ws.Range("K136").Value = 10959.1032
$ws.Range("L136").Value = 14820
$ws.Range("M136").Value = -8409.1032
$ws.Range("N136").Value = -19920
